# 自动更新Excel文件
# For every data row (2..99) on the active sheet, decrement the
# "剩余" (remaining days) value in column E by 1.
# When the remaining-days counter reaches 0, it rolls over to 10
# (a new 10-day cycle) and the "开始时间" (start date, column F,
# stored as an 8-digit yyyymmdd number) is advanced by 10 days.
#
# Rows whose start date is not a valid 8-digit yyyymmdd number are
# left untouched (their data is considered malformed / out of scope).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $remainingCell = $ws.Cells.Item($r, 5)   # Column E - 剩余
    $startCell     = $ws.Cells.Item($r, 6)   # Column F - 开始时间

    $remaining = $remainingCell.Value()
    $startDate = $startCell.Value()

    if ($null -eq $remaining -or $null -eq $startDate) {
        continue
    }

    $startDateStr = [string][int64]$startDate

    # Only process rows with a well-formed 8-digit yyyymmdd date
    if ($startDateStr.Length -ne 8) {
        continue
    }

    $newRemaining = [int]$remaining - 1

    if ($newRemaining -le 0) {
        $newRemaining = 10

        $dt = [datetime]::ParseExact($startDateStr, "yyyyMMdd", $null)
        $dt = $dt.AddDays(10)

        $newStartDate = [int64]($dt.ToString("yyyyMMdd"))
        $startCell.Value = $newStartDate
    }

    $remainingCell.Value = $newRemaining
}
